# Auto-generated edits applying the Mateus_Profits.xlsx diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 7254.522
$ws.Range("I15").Value = 7254.522
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 21763.566
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -21594.566

$ws.Range("H19").Value = 4445.875
$ws.Range("I19").Value = 4428
$ws.Range("J19").Value = 4448.4287
$ws.Range("K19").Value = 4428
$ws.Range("L19").Value = 4448.4287
$ws.Range("M19").Value = -4253
$ws.Range("N19").Value = -4798.4287

$ws.Range("H28").Value = 589.53845
$ws.Range("I28").Value = 715.5
$ws.Range("J28").Value = 388
$ws.Range("K28").Value = 715.5
$ws.Range("L28").Value = 388
$ws.Range("M28").Value = -230.5
$ws.Range("N28").Value = -1358

$ws.Range("H33").Value = 320
$ws.Range("I33").Value = 141.33333
$ws.Range("J33").Value = 3000
$ws.Range("K33").Value = 141.33333
$ws.Range("L33").Value = 3000
$ws.Range("M33").Value = 87.66667000000001
$ws.Range("N33").Value = -3458

$ws.Range("H55").Value = 552.875
$ws.Range("I55").Value = 484
$ws.Range("J55").Value = 667.6667
$ws.Range("K55").Value = 484
$ws.Range("L55").Value = 667.6667
$ws.Range("M55").Value = -270
$ws.Range("N55").Value = -1095.6667

$ws.Range("H58").Value = 74.333336
$ws.Range("I58").Value = 74.333336
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 223.000008
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -73.00000800000001
$ws.Range("N58").ClearContents()

$ws.Range("H138").Value = 2429.5386
$ws.Range("I138").Value = 2178
$ws.Range("J138").Value = 2475.2727
$ws.Range("K138").Value = 6534
$ws.Range("L138").Value = 7425.8181
$ws.Range("M138").Value = -1394
$ws.Range("N138").Value = -17705.8181

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 178.15384
$ws.Range("I5").Value = 121.2
$ws.Range("J5").Value = 368
$ws.Range("K5").Value = 121.2
$ws.Range("L5").Value = 368
$ws.Range("M5").Value = -9.200000000000003
$ws.Range("N5").Value = -592

$ws.Range("H32").Value = 6030.961
$ws.Range("I32").Value = 5868.959
$ws.Range("J32").Value = 10000
$ws.Range("K32").Value = 5868.959
$ws.Range("L32").Value = 10000
$ws.Range("M32").Value = -5581.959
$ws.Range("N32").Value = -10574

$ws.Range("H61").Value = 6801.7856
$ws.Range("I61").Value = 5194.826
$ws.Range("J61").Value = 14193.8
$ws.Range("K61").Value = 5194.826
$ws.Range("L61").Value = 14193.8
$ws.Range("M61").Value = -4982.826
$ws.Range("N61").Value = -14617.8

$ws.Range("H110").Value = 5829.65
$ws.Range("I110").Value = 4199.6924
$ws.Range("J110").Value = 8856.714
$ws.Range("K110").Value = 4199.6924
$ws.Range("L110").Value = 8856.714
$ws.Range("M110").Value = -2154.6924
$ws.Range("N110").Value = -12946.714

$ws.Range("H136").Value = 6801.7856
$ws.Range("I136").Value = 5194.826
$ws.Range("J136").Value = 14193.8
$ws.Range("K136").Value = 15584.478
$ws.Range("L136").Value = 42581.39999999999
$ws.Range("M136").Value = -13034.478
$ws.Range("N136").Value = -47681.39999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 178.15384
$ws.Range("I4").Value = 121.2
$ws.Range("J4").Value = 368
$ws.Range("K4").Value = 121.2
$ws.Range("L4").Value = 368
$ws.Range("M4").Value = -6.200000000000003
$ws.Range("N4").Value = -598

$ws.Range("H64").Value = 673.5714
$ws.Range("I64").Value = 625
$ws.Range("J64").Value = 693
$ws.Range("K64").Value = 625
$ws.Range("L64").Value = 693
$ws.Range("M64").Value = -400
$ws.Range("N64").Value = -1143

$ws.Range("H67").Value = 673.5714
$ws.Range("I67").Value = 625
$ws.Range("J67").Value = 693
$ws.Range("K67").Value = 625
$ws.Range("L67").Value = 693
$ws.Range("M67").Value = 155
$ws.Range("N67").Value = -2253

$ws.Range("H80").Value = 754
$ws.Range("I80").Value = 500
$ws.Range("J80").Value = 838.6667
$ws.Range("K80").Value = 500
$ws.Range("L80").Value = 838.6667
$ws.Range("M80").Value = 498
$ws.Range("N80").Value = -2834.6667

$ws.Range("H83").Value = 754
$ws.Range("I83").Value = 500
$ws.Range("J83").Value = 838.6667
$ws.Range("K83").Value = 2500
$ws.Range("L83").Value = 4193.3335
$ws.Range("M83").Value = 2492
$ws.Range("N83").Value = -14177.3335

$ws.Range("H86").Value = 2001.4166
$ws.Range("I86").Value = 1774.2727
$ws.Range("J86").Value = 4500
$ws.Range("K86").Value = 1774.2727
$ws.Range("L86").Value = 4500
$ws.Range("M86").Value = -651.2727
$ws.Range("N86").Value = -6746

$ws.Range("H89").Value = 2001.4166
$ws.Range("I89").Value = 1774.2727
$ws.Range("J89").Value = 4500
$ws.Range("K89").Value = 8871.363499999999
$ws.Range("L89").Value = 22500
$ws.Range("M89").Value = -3255.363499999999
$ws.Range("N89").Value = -33732

$ws.Range("H97").Value = 6378.8335
$ws.Range("I97").Value = 6378.8335
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 6378.8335
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -5387.8335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H8").Value = 2750
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 2750
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 2750
$ws.Range("N8").Value = -3030

$ws.Range("H105").Value = 1949.2667
$ws.Range("I105").Value = 1989.5714
$ws.Range("J105").Value = 1385
$ws.Range("K105").Value = 1989.5714
$ws.Range("L105").Value = 1385
$ws.Range("M105").Value = -242.5714
$ws.Range("N105").Value = -4879

$ws.Range("H107").Value = 531.7083
$ws.Range("I107").Value = 460
$ws.Range("J107").Value = 632.1
$ws.Range("K107").Value = 460
$ws.Range("L107").Value = 632.1
$ws.Range("M107").Value = 1460
$ws.Range("N107").Value = -4472.1

$ws.Range("H129").Value = 69498
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 69498
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 69498
$ws.Range("N129").Value = -79498

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 877.5
$ws.Range("I22").Value = 1133
$ws.Range("J22").Value = 111
$ws.Range("K22").Value = 3399
$ws.Range("L22").Value = 333
$ws.Range("M22").Value = -3230
$ws.Range("N22").Value = -671

$ws.Range("H27").Value = 877.5
$ws.Range("I27").Value = 1133
$ws.Range("J27").Value = 111
$ws.Range("K27").Value = 3399
$ws.Range("L27").Value = 333
$ws.Range("M27").Value = -3297
$ws.Range("N27").Value = -537

$ws.Range("H40").Value = 175.75
$ws.Range("I40").Value = 68.833336
$ws.Range("J40").Value = 239.9
$ws.Range("K40").Value = 275.333344
$ws.Range("L40").Value = 959.6
$ws.Range("M40").Value = -206.333344
$ws.Range("N40").Value = -1097.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2920.3044
$ws.Range("I97").Value = 537.2222
$ws.Range("J97").Value = 11499.4
$ws.Range("K97").Value = 537.2222
$ws.Range("L97").Value = 11499.4
$ws.Range("M97").Value = -41.22220000000004
$ws.Range("N97").Value = -12491.4

$ws.Range("H107").Value = 862437.5
$ws.Range("I107").Value = 1567712.5
$ws.Range("J107").Value = 434.66666
$ws.Range("K107").Value = 1567712.5
$ws.Range("L107").Value = 434.66666
$ws.Range("M107").Value = -1565792.5
$ws.Range("N107").Value = -4274.66666

$ws.Range("H129").Value = 64749.75
$ws.Range("I129").Value = 60000
$ws.Range("J129").Value = 66333
$ws.Range("K129").Value = 60000
$ws.Range("L129").Value = 66333
$ws.Range("M129").Value = -55000
$ws.Range("N129").Value = -76333

$ws.Range("H132").Value = 3750.875
$ws.Range("I132").Value = 3048.5386
$ws.Range("J132").Value = 6794.3335
$ws.Range("K132").Value = 9145.6158
$ws.Range("L132").Value = 20383.0005
$ws.Range("M132").Value = -6615.6158
$ws.Range("N132").Value = -25443.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4654.778
$ws.Range("I22").Value = 5066.3335
$ws.Range("J22").Value = 3831.6667
$ws.Range("K22").Value = 5066.3335
$ws.Range("L22").Value = 3831.6667
$ws.Range("M22").Value = -4771.3335
$ws.Range("N22").Value = -4421.6667

$ws.Range("H27").Value = 4654.778
$ws.Range("I27").Value = 5066.3335
$ws.Range("J27").Value = 3831.6667
$ws.Range("K27").Value = 5066.3335
$ws.Range("L27").Value = 3831.6667
$ws.Range("M27").Value = -4959.3335
$ws.Range("N27").Value = -4045.6667

$ws.Range("H132").Value = 12374.682
$ws.Range("I132").Value = 13014
$ws.Range("J132").Value = 9497.75
$ws.Range("K132").Value = 39042
$ws.Range("L132").Value = 28493.25
$ws.Range("M132").Value = -36512
$ws.Range("N132").Value = -33553.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3506.1924
$ws.Range("I132").Value = 3014.5715
$ws.Range("J132").Value = 5571
$ws.Range("K132").Value = 9043.7145
$ws.Range("L132").Value = 16713
$ws.Range("M132").Value = -6513.7145
$ws.Range("N132").Value = -21773

$ws.Range("H136").Value = 3216.7036
$ws.Range("I136").Value = 2683.182
$ws.Range("J136").Value = 5564.2
$ws.Range("K136").Value = 8049.545999999999
$ws.Range("L136").Value = 16692.6
$ws.Range("M136").Value = -5499.545999999999
$ws.Range("N136").Value = -21792.6
